$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the trial-length header values (row 1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update CON row (row 2) values
$ws.Range("B2").Value = 37.31198412568245
$ws.Range("C2").Value = 45.407415285381461
$ws.Range("D2").Value = 40.905546271713845
$ws.Range("E2").Value = 48.242823049341752

# Update STR row (row 3) values
$ws.Range("B3").Value = 43.216688876332171
$ws.Range("C3").Value = 55.236951175552917
$ws.Range("D3").Value = 54.251121620335695
$ws.Range("E3").Value = 47.760437452000765

# Update the active selection to reflect the reduced range
$ws.Range("B1:E3").Select()
